$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D and E) for the reduced-emissions-for-cars scenarios,
# shifting the existing "Emissions for boats" column (and everything after it) right.
$ws.Range("D1:E1").EntireColumn.Insert()

$ws.Range("D1").Value = "Reduced Emissions for cars (5%)"
$ws.Range("E1").Value = "Reduced Emissions for cars (10%)"

$ws.Range("D2").Formula = "=0.95*C2"
$ws.Range("E2").Formula = "=0.9*C2"
$ws.Range("D3:D7").Formula = "=0.95*C3"
$ws.Range("E3:E7").Formula = "=0.9*C3"

$ws.Rows("1").RowHeight = 58

$ws.Range("E10").Select() | Out-Null
